$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "64.096.92"
$ws.Range("E2").Value = "  -0.71%  "

# Row 3
$ws.Range("D3").Value = "3.319.45"
$ws.Range("E3").Value = "  -1.45%  "

# Row 4
$ws.Range("E4").Value = "  +0.04%  "

# Row 5
$styleD5 = $ws.Range("D5").Style
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "550.73"
$ws.Range("D5").Style = $styleD5
$ws.Range("E5").Value = "  -0.93%  "

# Row 6
$styleD6 = $ws.Range("D6").Style
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "172.46"
$ws.Range("D6").Style = $styleD6
$ws.Range("E6").Value = "  -2.19%  "

# Row 7
$styleD7 = $ws.Range("D7").Style
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.624"
$ws.Range("D7").Style = $styleD7
$ws.Range("E7").Value = "  +1.05%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("D9").Value = "3.312.53"
$ws.Range("E9").Value = "  -1.48%  "

# Row 10
$styleD10 = $ws.Range("D10").Style
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.170"
$ws.Range("D10").Style = $styleD10
$ws.Range("E10").Value = "  +5.52%  "

# Row 11
$styleD11 = $ws.Range("D11").Style
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.636"
$ws.Range("D11").Style = $styleD11
$ws.Range("E11").Value = "  +1.29%  "

# Row 12
$styleD12 = $ws.Range("D12").Style
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "53.06"
$ws.Range("D12").Style = $styleD12
$ws.Range("E12").Value = "  -3.37%  "

# Row 13
$styleD13 = $ws.Range("D13").Style
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000275"
$ws.Range("D13").Style = $styleD13
$ws.Range("E13").Value = "  +1.06%  "

# Row 14
$styleD14 = $ws.Range("D14").Style
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "9.03"
$ws.Range("D14").Style = $styleD14
$ws.Range("E14").Value = "  -0.39%  "

# Row 15
$ws.Range("D15").Value = "3.854.74"
$ws.Range("E15").Value = "  -1.36%  "

# Row 16
$ws.Range("E16").Value = "  +2.12%  "

# Row 17
$styleD17 = $ws.Range("D17").Style
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.05"
$ws.Range("D17").Style = $styleD17
$ws.Range("E17").Value = "  -1.94%  "

# Row 18
$ws.Range("D18").Value = "3.314.32"
$ws.Range("E18").Value = "  -1.79%  "

# Row 19
$ws.Range("D19").Value = "63.985.24"
$ws.Range("E19").Value = "  -0.76%  "

# Row 20
$styleD20 = $ws.Range("D20").Style
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.67"
$ws.Range("D20").Style = $styleD20
$ws.Range("E20").Value = "  -1.46%  "

# Row 21
$styleD21 = $ws.Range("D21").Style
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "0.980"
$ws.Range("D21").Style = $styleD21
$ws.Range("E21").Value = "  -0.20%  "

# Row 22
$styleD22 = $ws.Range("D22").Style
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "446.33"
$ws.Range("D22").Style = $styleD22
$ws.Range("E22").Value = "  +3.41%  "

# Row 23
$styleD23 = $ws.Range("D23").Style
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.97"
$ws.Range("D23").Style = $styleD23
$ws.Range("E23").Value = "  +0.57%  "

# Row 24
$ws.Range("E24").Value = "  -2.03%  "

# Row 25
$styleD25 = $ws.Range("D25").Style
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "86.55"
$ws.Range("D25").Style = $styleD25
$ws.Range("E25").Value = "  +2.76%  "

# Row 26
$styleD26 = $ws.Range("D26").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "13.69"
$ws.Range("D26").Style = $styleD26
$ws.Range("E26").Value = "  +3.14%  "

# Row 27
$ws.Range("E27").Value = "  +0.88%  "

# Row 28
$ws.Range("E28").Value = "  -1.95%  "

# Row 29
$styleD29 = $ws.Range("D29").Style
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "8.54"
$ws.Range("D29").Style = $styleD29
$ws.Range("E29").Value = "  -2.80%  "

# Row 30
$styleD30 = $ws.Range("D30").Style
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "30.65"
$ws.Range("D30").Style = $styleD30
$ws.Range("E30").Value = "  +3.04%  "

# Row 31
$styleD31 = $ws.Range("D31").Style
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "6.48"
$ws.Range("D31").Style = $styleD31
$ws.Range("E31").Value = "  -2.91%  "

# Row 32
$ws.Range("B32").Value = "OKB"
$ws.Range("C32").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$styleD32 = $ws.Range("D32").Style
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "62.16"
$ws.Range("D32").Style = $styleD32
$ws.Range("E32").Value = "  +6.26%  "

# Row 33
$ws.Range("B33").Value = "Cosmos"
$ws.Range("C33").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$styleD33 = $ws.Range("D33").Style
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "11.32"
$ws.Range("D33").Style = $styleD33
$ws.Range("E33").Value = "  -1.49%  "

# Row 34
$styleD34 = $ws.Range("D34").Style
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "568.55"
$ws.Range("D34").Style = $styleD34
$ws.Range("E34").Value = "  -0.85%  "

# Row 35
$ws.Range("E35").Value = "  -1.80%  "

# Row 37
$ws.Range("E37").Value = "  +0.70%  "

# Row 38
$ws.Range("E38").Value = "  -1.17%  "

# Row 39
$styleD39 = $ws.Range("D39").Style
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "35.06"
$ws.Range("D39").Style = $styleD39
$ws.Range("E39").Value = "  -2.00%  "

# Row 40
$ws.Range("E40").Value = "  -1.08%  "

# Row 41
$ws.Range("D41").Value = "0.0₃0726"
$ws.Range("E41").Value = "  -4.16%  "

# Row 42
$ws.Range("D42").Value = "3.051.98"
$ws.Range("E42").Value = "  -2.02%  "

# Row 43
$ws.Range("E43").Value = "  +0.43%  "

# Row 44
$ws.Range("E44").Value = "  -4.12%  "

# Row 45
$ws.Range("E45").Value = "  -3.12%  "

# Row 46
$ws.Range("E46").Value = "  +2.29%  "

# Row 47
$styleD47 = $ws.Range("D47").Style
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.42"
$ws.Range("D47").Style = $styleD47
$ws.Range("E47").Value = "  -1.72%  "

# Row 48
$ws.Range("B48").Value = "FirstDigitalUSD"
$ws.Range("C48").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$styleD48 = $ws.Range("D48").Style
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.00"
$ws.Range("D48").Style = $styleD48
$ws.Range("E48").Value = "  +0.06%  "

# Row 49
$ws.Range("B49").Value = "Monero"
$ws.Range("C49").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$styleD49 = $ws.Range("D49").Style
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "142.43"
$ws.Range("D49").Style = $styleD49
$ws.Range("E49").Value = "  +5.59%  "

# Row 50
$ws.Range("E50").Value = "  -3.10%  "

# Row 51
$styleD51 = $ws.Range("D51").Style
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.14"
$ws.Range("D51").Style = $styleD51
$ws.Range("E51").Value = "  -1.44%  "
